$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: new "No LogHole to WellTop" well list ---
$ws.Range("C1").Value = 'No LogHole to WellTop'
$ws.Range("C2").Value = 'EEF-1'
$ws.Range("C3").Value = 'EER-1'
$ws.Range("C4").Value = 'EET-1'
$ws.Range("C5").Value = 'EET-2'
$ws.Range("C6").Value = 'EHR-1'
$ws.Range("C7").Value = 'EWF-XA2S'
$ws.Range("C8").Value = 'A3S'
$ws.Range("C9").Value = 'EHT-XA5S'
$ws.Range("C10").Value = 'EWR-1'
$ws.Range("C11").Value = 'EWT-1'
$ws.Range("C13").Value = 'KEF-1'
$ws.Range("C14").Value = 'KET-1'
$ws.Range("C15").Value = 'KEFR-XE1A'
$ws.Range("C16").Value = 'KER-1'
$ws.Range("C17").Value = 'KFR-XE1S'
$ws.Range("C18").Value = 'KHT-XE5S'
$ws.Range("C19").Value = 'KWF-1'
$ws.Range("C20").Value = 'KWR-1'
$ws.Range("C21").Value = 'KWT-1'
$ws.Range("C22").Value = 'KHF-1'
$ws.Range("C23").Value = 'KHT-XE4S'
$ws.Range("C24").Value = 'E5S'
$ws.Range("C26").Value = 'LHR-1'
$ws.Range("C27").Value = 'LHR-1'

# Column C width
$ws.Columns("C").ColumnWidth = 17

# Apply the Arial 10pt font to the subset of C cells that use it
$cArialRange = $ws.Range("C5,C7,C14,C22,C23,C24,C26,C27")
$cArialRange.Font.Name = "Arial"
$cArialRange.Font.Size = 10

# --- Column F: formatted-but-empty helper cells (Arial 10pt) ---
$fRange = $ws.Range("F9,F10,F11,F12,F13,F14,F16,F17,F19,F20,F21,F22,F23,F24,F25,F26,F27,F28,F30,F31,F32,F33,F34,F35,F36,F37,F38,F39,F40,F41,F42,F43,F44,F45,F46")
$fRange.Font.Name = "Arial"
$fRange.Font.Size = 10

# --- Sort state metadata over F7:F46 (mirrors a prior sort on helper column F) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("F7:F46"))
$ws.Sort.SetRange($ws.Range("F7:F46"))
$ws.Sort.Apply()

# --- Selection & view state ---
$ws.Range("C2:C27").Select()
